# Daily attendance processing - reorder the "Recorded By" (column G) names.
#
# The values in column G are comma-separated lists of the people/accounts
# that recorded a session (e.g. "dnasr281@gmail.com, System"). This pass
# re-normalizes the ordering of those lists according to a fixed, case
# sensitive priority so the same accounts are always listed in the same
# relative order:
#   admin@admin.com > System > system > dnasr281@gmail.com > backup@backdoor.com
# Any name not in the known list keeps its relative position after the
# known ones (stable sort), and cells without a comma are left untouched.
#
# NOTE: this engine's functions do not get their own variable scope (loop
# counters leak into the caller), and string comparisons/hashtable lookups
# are case-insensitive, so helper functions below use distinctively-named
# variables and compare character codes by hand to stay case sensitive.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CaseSensitiveEquals($cse_s1, $cse_s2) {
    if ($cse_s1.Length -ne $cse_s2.Length) { return $false }
    for ($cse_k = 0; $cse_k -lt $cse_s1.Length; $cse_k++) {
        $cse_c1 = [int][char]$cse_s1.Substring($cse_k, 1)
        $cse_c2 = [int][char]$cse_s2.Substring($cse_k, 1)
        if ($cse_c1 -ne $cse_c2) { return $false }
    }
    return $true
}

function Get-NamePriority($gnp_name) {
    if (CaseSensitiveEquals $gnp_name "admin@admin.com") { return 0 }
    if (CaseSensitiveEquals $gnp_name "System") { return 1 }
    if (CaseSensitiveEquals $gnp_name "system") { return 2 }
    if (CaseSensitiveEquals $gnp_name "dnasr281@gmail.com") { return 3 }
    if (CaseSensitiveEquals $gnp_name "backup@backdoor.com") { return 4 }
    return 100
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($rowNum = 1; $rowNum -le $lastRow; $rowNum++) {
    $gCell = $ws.Cells.Item($rowNum, 7)   # Column G = "Recorded By"
    $origText = $gCell.Text

    if ([string]::IsNullOrEmpty($origText)) { continue }
    if ($origText -notmatch ",") { continue }

    $nameParts = $origText -split ",\s*" | ForEach-Object { $_.Trim() }

    $entries = @()
    for ($partIdx = 0; $partIdx -lt $nameParts.Count; $partIdx++) {
        $entries += [PSCustomObject]@{
            Name = $nameParts[$partIdx]
            Pri  = Get-NamePriority $nameParts[$partIdx]
            Orig = $partIdx
        }
    }

    # NOTE: "-Property Pri, Orig" (comma list) is mis-parsed by this engine;
    # an explicit array is required for a correct multi-key sort.
    $sortedEntries = $entries | Sort-Object -Property @('Pri', 'Orig')
    $rebuiltText = ($sortedEntries | ForEach-Object { $_.Name }) -join ", "

    # NOTE: "-ne"/"-eq" string comparisons are case-insensitive in this
    # engine, so always (re)write the value instead of trying to skip
    # cells that are "unchanged" - writing the same text back is harmless.
    $gCell.Value = $rebuiltText
}
